$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "calculate_area" test-plan row: append the actual calculated
# perimeter / area values underneath their labels (new test results).
$ws.Range("G14").Value = "Calculated Perimeter" + [char]10 + "17"
$ws.Range("G13").Value = "Calculated Area" + [char]10 + "12.497499749949988"

# Scroll the sheet down so row 9 is at the top and select G13, matching
# where the author was working when they added the test results.
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
$ws.Range("G13").Select()
